$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2,4).Value = "42.834.02"
$ws.Cells.Item(2,5).Value = "  -2.24%  "
$ws.Cells.Item(3,4).Value = "2.572.81"
$ws.Cells.Item(3,5).Value = "  -0.33%  "
$ws.Cells.Item(4,5).Value = "  -0.11%  "
$ws.Cells.Item(5,4).Value = "'302.18"
$ws.Cells.Item(5,5).Value = "  +0.17%  "
$ws.Cells.Item(6,4).Value = "'96.73"
$ws.Cells.Item(6,5).Value = "  +0.07%  "
$ws.Cells.Item(7,5).Value = "  -1.11%  "
$ws.Cells.Item(8,5).Value = "  -0.03%  "
$ws.Cells.Item(9,4).Value = "'0.549"
$ws.Cells.Item(9,5).Value = "  -2.23%  "
$ws.Cells.Item(10,5).Value = "  -1.28%  "
$ws.Cells.Item(11,4).Value = "'0.0812"
$ws.Cells.Item(11,5).Value = "  -0.93%  "
$ws.Cells.Item(12,4).Value = "'7.64"
$ws.Cells.Item(12,5).Value = "  -2.32%  "
$ws.Cells.Item(13,5).Value = "  +6.26%  "
$ws.Cells.Item(14,4).Value = "2.612.24"
$ws.Cells.Item(14,5).Value = "  +0.94%  "
$ws.Cells.Item(15,4).Value = "'0.886"
$ws.Cells.Item(15,5).Value = "  -0.95%  "
$ws.Cells.Item(16,4).Value = "'14.34"
$ws.Cells.Item(16,5).Value = "  -0.26%  "
$ws.Cells.Item(17,4).Value = "42.895.03"
$ws.Cells.Item(17,5).Value = "  -2.05%  "
$ws.Cells.Item(18,4).Value = "0.0₃0997"
$ws.Cells.Item(18,5).Value = "  +1.31%  "
$ws.Cells.Item(19,4).Value = "'12.90"
$ws.Cells.Item(19,5).Value = "  +3.56%  "
$ws.Cells.Item(20,4).Value = "'6.64"
$ws.Cells.Item(20,5).Value = "  -1.16%  "
$ws.Cells.Item(21,4).Value = "'72.10"
$ws.Cells.Item(21,5).Value = "  -2.06%  "
$ws.Cells.Item(22,4).Value = "'254.41"
$ws.Cells.Item(22,5).Value = "  -4.36%  "
$ws.Cells.Item(23,5).Value = "  +0.37%  "
$ws.Cells.Item(24,5).Value = "  -5.47%  "
$ws.Cells.Item(25,4).Value = "'29.03"
$ws.Cells.Item(25,5).Value = "  -1.40%  "
$ws.Cells.Item(26,5).Value = "  -0.18%  "
$ws.Cells.Item(27,4).Value = "'10.28"
$ws.Cells.Item(27,5).Value = "  +0.15%  "
$ws.Cells.Item(28,4).Value = "'37.63"
$ws.Cells.Item(28,5).Value = "  -1.67%  "
$ws.Cells.Item(29,5).Value = "  -5.65%  "
$ws.Cells.Item(30,4).Value = "'6.01"
$ws.Cells.Item(30,5).Value = "  -2.78%  "
$ws.Cells.Item(31,4).Value = "'154.69"
$ws.Cells.Item(31,5).Value = "  +1.40%  "
$ws.Cells.Item(32,5).Value = "  -2.12%  "
$ws.Cells.Item(33,4).Value = "'3.40"
$ws.Cells.Item(33,5).Value = "  -4.99%  "
$ws.Cells.Item(34,5).Value = "  -1.74%  "
$ws.Cells.Item(35,4).Value = "'0.0802"
$ws.Cells.Item(35,5).Value = "  -1.89%  "
$ws.Cells.Item(36,4).Value = "'18.29"
$ws.Cells.Item(36,5).Value = "  +7.83%  "
$ws.Cells.Item(37,5).Value = "  -3.23%  "
$ws.Cells.Item(38,5).Value = "  -0.81%  "
$ws.Cells.Item(39,4).Value = "'22.94"
$ws.Cells.Item(39,5).Value = "  -5.70%  "
$ws.Cells.Item(40,4).Value = "'2.10"
$ws.Cells.Item(40,5).Value = "  +29.11%  "
$ws.Cells.Item(41,5).Value = "  -4.95%  "
$ws.Cells.Item(42,2).Value = "VeChain"
$ws.Cells.Item(42,3).Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Cells.Item(42,4).Value = "'0.0311"
$ws.Cells.Item(42,5).Value = "  -1.84%  "
$ws.Cells.Item(43,2).Value = "RenderToken"
$ws.Cells.Item(43,3).Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Cells.Item(43,4).Value = "'3.88"
$ws.Cells.Item(43,5).Value = "  -0.10%  "
$ws.Cells.Item(44,4).Value = "2.080.91"
$ws.Cells.Item(44,5).Value = "  +2.09%  "
$ws.Cells.Item(45,4).Value = "'0.998"
$ws.Cells.Item(45,5).Value = "  +0.01%  "
$ws.Cells.Item(46,5).Value = "  +0.25%  "
$ws.Cells.Item(47,4).Value = "'85.17"
$ws.Cells.Item(47,5).Value = "  -3.26%  "
$ws.Cells.Item(48,4).Value = "'76.04"
$ws.Cells.Item(48,5).Value = "  +9.52%  "
$ws.Cells.Item(49,4).Value = "'106.36"
$ws.Cells.Item(49,5).Value = "  +0.74%  "
$ws.Cells.Item(50,4).Value = "2.821.64"
$ws.Cells.Item(50,5).Value = "  -0.60%  "
$ws.Cells.Item(51,5).Value = "  -0.32%  "
